$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.203.42'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '1.828.98'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.15'
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6101'
$ws.Range("E6").Value = '  -3.87%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07094'
$ws.Range("E8").Value = '  -5.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.92'
$ws.Range("E10").Value = '  -4.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07638'
$ws.Range("E11").Value = '  -1.37%  '
$ws.Range("D12").Value = '1.836.78'
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6385'
$ws.Range("E14").Value = '  -6.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009959'
$ws.Range("E15").Value = '  -2.76%  '
$ws.Range("D16").Value = '2.071.19'
$ws.Range("E16").Value = '  -1.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.36'
$ws.Range("E17").Value = '  -3.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.963'
$ws.Range("E18").Value = '  -4.65%  '
$ws.Range("D19").Value = '29.205.33'
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '229.11'
$ws.Range("E20").Value = '  -0.46%  '
$ws.Range("E21").Value = '  -4.23%  '
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("E23").Value = '  -4.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.52'
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("E26").Value = '  -4.61%  '
$ws.Range("E27").Value = '  -4.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.73'
$ws.Range("E28").Value = '  -4.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06759'
$ws.Range("E29").Value = '  +3.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.487'
$ws.Range("E30").Value = '  +3.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.461'
$ws.Range("E31").Value = '  -1.85%  '
$ws.Range("E32").Value = '  -5.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.835'
$ws.Range("E33").Value = '  -5.28%  '
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.741'
$ws.Range("E35").Value = '  -5.49%  '
$ws.Range("E37").Value = '  -0.86%  '
$ws.Range("D38").Value = '1.235.49'
$ws.Range("E38").Value = '  -1.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.760'
$ws.Range("E40").Value = '  -4.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.589'
$ws.Range("E41").Value = '  -2.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9278'
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '1.983.63'
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.02'
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.66'
$ws.Range("E46").Value = '  -2.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000115'
$ws.Range("E47").Value = '  -3.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.630'
$ws.Range("E48").Value = '  -5.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.575'
$ws.Range("E49").Value = '  -5.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.544'
$ws.Range("E50").Value = '  -7.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1085'
$ws.Range("E51").Value = '  -5.56%  '
